$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shipment-tracking numbers (shared-string text, not numbers) landed
# on the cheetah/baseinit test rows -- P2:P5 get new ShipmentTracking
# values. Force text storage (so they stay shared strings like the
# originals) and then restore the default "Normal" style so no visible
# formatting changes, matching the source data's unstyled cells.
$ws.Range("P2:P5").NumberFormat = "@"

$ws.Range("P2").Value = "320018656023"
$ws.Range("P3").Value = "320018646044"
$ws.Range("P4").Value = "320018646077"
$ws.Range("P5").Value = "320018624977"

$ws.Range("P2:P5").Style = "Normal"
